$d = $word.ActiveDocument
$bullet = [char]0x2022

# --- Change 1: Collapse the three "CORE COMPETENCIES" detail paragraphs into a
# single summary line ---------------------------------------------------------
$pCore = $d.Paragraphs.Item(6)
$rCore = $d.Range($pCore.Range.Start, $pCore.Range.End - 1)
$rCore.Text = "Product Marketing Core " + $bullet + " Research & Analytics " + $bullet + " Communication & Technology"

# Remove the two now-redundant paragraphs that followed it (their detailed text
# moves to the new TECHNICAL SKILLS section below).
$pResearch = $d.Paragraphs.Item(7)
$pComm = $d.Paragraphs.Item(8)
$delRange = $d.Range($pResearch.Range.Start, $pComm.Range.End)
$delRange.Delete()

# --- Change 2: Add a new "TECHNICAL SKILLS" section with the detailed
# competency text, placed just before the closing LinkedIn/Site paragraph ----
$pLastBullet = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$pLastBullet.Range.InsertParagraphAfter()

$pHeading = $d.Paragraphs.Item($pLastBullet.Index + 1)
$pHeading.Range.Text = "TECHNICAL SKILLS"
$pHeading.Style = "Heading2"

$pHeading.Range.InsertParagraphAfter()
$pProduct = $d.Paragraphs.Item($pHeading.Index + 1)
$pProduct.Style = "Normal"
$pProduct.Range.Text = "PRODUCT MARKETING CORE Market Intelligence & Competitive Analysis; Product Positioning & Messaging Development; Go-to-Market Strategy & Product Launch Management; Customer Segmentation & Buyer Persona Development"

$pProduct.Range.InsertParagraphAfter()
$pResearch2 = $d.Paragraphs.Item($pProduct.Index + 1)
$pResearch2.Style = "Normal"
$pResearch2.Range.Text = "RESEARCH & ANALYTICS Survey Methodology & Customer Insights; Market Research Design & Implementation; Competitive Intelligence & SWOT Analysis; A/B Testing & Conversion Optimization"

$pResearch2.Range.InsertParagraphAfter()
$pComm2 = $d.Paragraphs.Item($pResearch2.Index + 1)
$pComm2.Style = "Normal"
$pComm2.Range.Text = "COMMUNICATION & TECHNOLOGY Strategic Messaging & Narrative Development; Technical Concept Translation for Business Audiences; Data Visualization & Reporting (Tableau, PowerBI, d3.js); Client Relationship Management & Business Development"
